# "Add files via upload" — Salary/W5 Salaries and Tasks.xlsx
#
# The only functional content change in this commit is the task note in the
# "Tasks to complete next week" column: the placeholder "tba" (row 19,
# column B) is replaced with the real task description. Everything else in
# the underlying XML diff (fileVersion/rupBuild, absPath, revisionPtr GUID,
# window geometry, cellXfs ordering, x14ac:dyDescent bookkeeping, etc.) is
# Excel-build resave noise, not an intentional edit, so we leave it alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "tba" placeholder to the real task description.
$ws.Range("B19").Value = 'Pick "perfect" idea and prototype it'

# Match the cursor/selection position left behind by the edit.
$ws.Range("B20").Select()
